$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.587.55'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').Value = '1.666.98'
$ws.Range('E3').Value = '  -3.44%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'215.24"
$ws.Range('E5').Value = '  -1.61%  '
$ws.Range('E6').Value = '  -1.85%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = "'23.58"
$ws.Range('E8').Value = '  -2.17%  '
$ws.Range('E9').Value = '  -1.23%  '
$ws.Range('E10').Value = '  -1.93%  '
$ws.Range('E11').Value = '  -2.46%  '
$ws.Range('D12').Value = '1.902.82'
$ws.Range('E12').Value = '  -3.41%  '
$ws.Range('D13').Value = '1.650.25'
$ws.Range('E13').Value = '  -4.30%  '
$ws.Range('E14').Value = '  -2.65%  '
$ws.Range('E15').Value = '  -2.24%  '
$ws.Range('E16').Value = '  -2.36%  '
$ws.Range('D17').Value = "'250.91"
$ws.Range('E17').Value = '  +2.67%  '
$ws.Range('D18').Value = '27.601.75'
$ws.Range('E18').Value = '  -1.13%  '
$ws.Range('E19').Value = '  -3.19%  '
$ws.Range('E20').Value = '  -4.33%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('E22').Value = '  -2.96%  '
$ws.Range('E23').Value = '  -4.69%  '
$ws.Range('E24').Value = '  -5.84%  '
$ws.Range('D25').Value = "'146.61"
$ws.Range('E25').Value = '  -1.85%  '
$ws.Range('E26').Value = '  -1.42%  '
$ws.Range('E27').Value = '  -5.08%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = "'0.112"
$ws.Range('E28').Value = '  -2.27%  '
$ws.Range('B29').Value = 'BinanceUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D29').Value = "'1.00"
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').Value = '  +4.01%  '
$ws.Range('D31').Value = "'0.0509"
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('D32').Value = "'3.35"
$ws.Range('E32').Value = '  -2.75%  '
$ws.Range('D33').Value = '1.471.98'
$ws.Range('E33').Value = '  -1.22%  '
$ws.Range('E34').Value = '  -5.39%  '
$ws.Range('E35').Value = '  -5.34%  '
$ws.Range('D36').Value = "'0.944"
$ws.Range('E36').Value = '  -1.77%  '
$ws.Range('E37').Value = '  -0.98%  '
$ws.Range('E38').Value = '  -5.72%  '
$ws.Range('E39').Value = '  -2.24%  '
$ws.Range('D40').Value = "'69.71"
$ws.Range('D41').Value = "'1.03"
$ws.Range('E41').Value = '  -3.71%  '
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('E43').Value = '  -7.03%  '
$ws.Range('D44').Value = '1.810.68'
$ws.Range('E44').Value = '  -3.37%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = "'0.791"
$ws.Range('E45').Value = '  -0.30%  '
$ws.Range('B46').Value = 'MXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D46').Value = "'2.20"
$ws.Range('E46').Value = '  -3.74%  '
$ws.Range('D47').Value = "'1.71"
$ws.Range('E47').Value = '  -2.08%  '
$ws.Range('D48').Value = "'89.51"
$ws.Range('E48').Value = '  -1.88%  '
$ws.Range('E49').Value = '  -2.18%  '
$ws.Range('D50').Value = "'42.08"
$ws.Range('E50').Value = '  +15.50%  '
$ws.Range('E51').Value = '  -3.42%  '
